# Generate Report for Handback
#
# This mirrors the "handback" step of the localization-status report: once a
# target (translated) file has been produced, the report is updated with:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (on the Overview sheet's per-locale
#     columns and on each locale sheet's Status column)
#   - The "Latest Target File" column gets the handed-back source file name
#     (hyperlinked, just like column A)
#   - The "Latest Handback File" column gets the xliff file name that was
#     produced for that locale
#   - The "Latest Handback DateTime" column gets the timestamp of the handback

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/"

$file1 = "ab6cf907-ef68-4f0f-9ae1-a47026b937b8.md"
$file2 = "d2f17490-2dad-4f74-9f05-6ee478d42efb.md"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn (col E) and de-de (col F) status cells
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("I2").Value = $file1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($ghBase + $file1), "", "", $file1) | Out-Null
$wsZh.Range("J2").Value = "ab6cf907-ef68-4f0f-9ae1-a47026b937b8.83c057c332cad23ee05fca9d24080b06bc355d72.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 04:58:45"

$wsZh.Range("I3").Value = $file2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($ghBase + $file2), "", "", $file2) | Out-Null
$wsZh.Range("J3").Value = "d2f17490-2dad-4f74-9f05-6ee478d42efb.c94ed268f6c6903dd1d3c94885ba664e6db4a140.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-26 04:58:45"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("I2").Value = $file1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($ghBase + $file1), "", "", $file1) | Out-Null
$wsDe.Range("J2").Value = "ab6cf907-ef68-4f0f-9ae1-a47026b937b8.83c057c332cad23ee05fca9d24080b06bc355d72.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 04:58:52"

$wsDe.Range("I3").Value = $file2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($ghBase + $file2), "", "", $file2) | Out-Null
$wsDe.Range("J3").Value = "d2f17490-2dad-4f74-9f05-6ee478d42efb.c94ed268f6c6903dd1d3c94885ba664e6db4a140.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-26 04:58:52"

# ---------------------------------------------------------------------------
# Column widths - widened to fit the newly-filled-in longer text
# (internal stored width = ColumnWidth + 5/6; pick ColumnWidth so stored
# width lands on the target)
# ---------------------------------------------------------------------------
$wideStatus = 30.0 - (5/6)      # -> stored width ~30 (status text grew)
$wideFile   = 40.0 - (5/6)      # -> stored width 40 (file-name columns)

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatus
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatus

$wsZh.Columns.Item(3).ColumnWidth = $wideStatus
$wsZh.Columns.Item(9).ColumnWidth = $wideFile
$wsZh.Columns.Item(10).ColumnWidth = $wideFile

$wsDe.Columns.Item(3).ColumnWidth = $wideStatus
$wsDe.Columns.Item(9).ColumnWidth = $wideFile
$wsDe.Columns.Item(10).ColumnWidth = $wideFile
